$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.578.79'
$ws.Range('E2').Value = '  -0.54%  '
$ws.Range('D3').Value = '2.286.41'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  +0.02%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '96.00'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +1.79%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '267.29'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -0.52%  '
$ws.Range('E7').Value = '  -1.62%  '
$ws.Range('E8').Value = '  -0.16%  '
$ws.Range('E9').Value = '  -2.04%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '45.85'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +0.38%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0930'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -0.40%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '7.79'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -3.13%  '
$ws.Range('E13').Value = '  -0.05%  '
$ws.Range('D14').Value = '2.626.36'
$ws.Range('E14').Value = '  -0.60%  '
$ws.Range('E15').Value = '  -1.27%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '0.846'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -0.12%  '
$ws.Range('D17').Value = '2.290.03'
$ws.Range('E17').Value = '  -1.32%  '
$ws.Range('D18').Value = '43.600.86'
$ws.Range('E18').Value = '  -0.40%  '
$ws.Range('E19').Value = '  +2.01%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '6.21'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -1.37%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '72.20'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +1.47%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '2.50'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +9.29%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '232.77'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -1.75%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '9.18'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -5.65%  '
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('E26').Value = '  +0.73%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '11.13'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -1.07%  '
$ws.Range('E28').Value = '  +2.55%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '40.56'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +3.92%  '
$ws.Range('E30').Value = '  +0.80%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '175.53'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +1.16%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '21.81'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -1.66%  '
$ws.Range('E33').Value = '  +0.22%  '
$ws.Range('E34').Value = '  -3.25%  '
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('E36').Value = '  -2.09%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.0354'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +1.51%  '
$ws.Range('E38').Value = '  -3.98%  '
$ws.Range('E39').Value = '  +0.76%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.236'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +0.18%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '12.24'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -0.45%  '
$ws.Range('B43').Value = 'MultiversX'
$ws.Range('C43').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '65.16'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +6.00%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '1.34'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +1.89%  '
$ws.Range('E45').Value = '  -0.76%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '5.21'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -4.64%  '
$ws.Range('E47').Value = '  -0.42%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '97.08'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -3.09%  '
$ws.Range('B49').Value = 'TrustWalletToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '1.19'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +0.19%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.432'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +0.45%  '
$ws.Range('D51').Value = '2.505.70'
$ws.Range('E51').Value = '  -0.64%  '
